# Generate Report for Handback
# The e6eb89bb-e7cd-4aeb-b010-5cbebbcb6dbd.md file has now been handed back
# (in sync with en-US), so update the status / datetime / error columns
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("K3").Value = "2016-11-14 06:28:18"
$wsZh.Range("P3").Value = ""

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("K3").Value = "2016-11-14 06:28:36"
$wsDe.Range("P3").Value = ""

# The Error Detail column on the language sheets used to need to be wide
# enough for the long error message; now that it is empty, shrink it back
# down to its auto-fit width (~13.75 characters).
$wsZh.Columns.Item(16).ColumnWidth = 12.8
$wsDe.Columns.Item(16).ColumnWidth = 12.8
